$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.417.86"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.818.13"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Formula = "'315.38"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").Formula = "'1.001"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  -4.11%  "
$ws.Range("D8").Formula = "'0.3959"
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("D9").Formula = "'0.08142"
$ws.Range("E9").Value = "  +7.26%  "
$ws.Range("D10").Formula = "'41.67"
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").Formula = "'1.108"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Formula = "'20.99"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").Formula = "'6.266"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("D14").Formula = "'1.001"
$ws.Range("E14").Value = "  +0.00%  "
$ws.Range("D15").Formula = "'7.503"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").Value = "1.818.47"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").Formula = "'0.00001144"
$ws.Range("E17").Value = "  +6.62%  "
$ws.Range("D18").Formula = "'92.61"
$ws.Range("E18").Value = "  +3.75%  "
$ws.Range("D19").Formula = "'0.06638"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").Formula = "'17.68"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Formula = "'1.001"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Formula = "'6.088"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "28.448.80"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Formula = "'11.28"
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("D25").Formula = "'2.267"
$ws.Range("E25").Value = "  +3.04%  "
$ws.Range("D26").Formula = "'21.14"
$ws.Range("E26").Value = "  +2.81%  "
$ws.Range("D27").Value = "2.029.26"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").Formula = "'155.20"
$ws.Range("E28").Value = "  -1.90%  "
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("D30").Formula = "'125.83"
$ws.Range("E30").Value = "  +1.67%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Formula = "'1.103"
$ws.Range("E32").Value = "  -1.54%  "
$ws.Range("D33").Formula = "'5.750"
$ws.Range("E33").Value = "  +1.73%  "
$ws.Range("D34").Formula = "'3.654"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("D35").Formula = "'0.07034"
$ws.Range("E35").Value = "  -4.50%  "
$ws.Range("D36").Formula = "'0.2228"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Formula = "'5.228"
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Formula = "'0.02326"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").Formula = "'8.818"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("D40").Formula = "'0.6267"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("D41").Formula = "'11.29"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").Formula = "'1.174"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Formula = "'1.401"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").Formula = "'13.51"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").Formula = "'3.739"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("D47").Formula = "'0.5911"
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("D48").Formula = "'124.90"
$ws.Range("D49").Formula = "'1.975"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("E51").Value = "  -0.07%  "
